$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced to text
# (matching the original workbook convention of storing all Price/Volume
# figures as text strings), without altering any other existing formatting.
$textCells = @("D5", "D6", "D8", "D9", "D11", "D15", "D16", "D18", "D20", "D23", "D25", "D26", "D29", "D34", "D37", "D38", "D40", "D42", "D43", "D44", "D45", "D46", "D48", "D49", "D51")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply the updated values
$ws.Range("D2").Value = "27.460.38"
$ws.Range("E2").Value = "  -0.98%  "
$ws.Range("D3").Value = "1.615.13"
$ws.Range("E3").Value = "  -1.82%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "211.04"
$ws.Range("E5").Value = "  -1.10%  "
$ws.Range("D6").Value = "0.526"
$ws.Range("E6").Value = "  -1.42%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "22.73"
$ws.Range("E8").Value = "  -1.95%  "
$ws.Range("D9").Value = "0.260"
$ws.Range("E9").Value = "  +0.40%  "
$ws.Range("E10").Value = "  -0.35%  "
$ws.Range("D11").Value = "0.0886"
$ws.Range("E11").Value = "  -0.44%  "
$ws.Range("D12").Value = "1.844.17"
$ws.Range("E12").Value = "  -1.80%  "
$ws.Range("D13").Value = "1.616.55"
$ws.Range("E13").Value = "  -1.83%  "
$ws.Range("E14").Value = "  -0.31%  "
$ws.Range("D15").Value = "0.548"
$ws.Range("E15").Value = "  -2.73%  "
$ws.Range("D16").Value = "64.87"
$ws.Range("E16").Value = "  +0.91%  "
$ws.Range("D17").Value = "27.447.97"
$ws.Range("E17").Value = "  -0.92%  "
$ws.Range("D18").Value = "231.00"
$ws.Range("E18").Value = "  -0.55%  "
$ws.Range("D19").Value = "0.0₃0718"
$ws.Range("E19").Value = "  -1.08%  "
$ws.Range("D20").Value = "7.51"
$ws.Range("E20").Value = "  -2.04%  "
$ws.Range("E21").Value = "  +0.14%  "
$ws.Range("E22").Value = "  -1.06%  "
$ws.Range("D23").Value = "10.15"
$ws.Range("E23").Value = "  +0.45%  "
$ws.Range("E24").Value = "  +5.63%  "
$ws.Range("D25").Value = "150.88"
$ws.Range("E25").Value = "  +0.67%  "
$ws.Range("D26").Value = "6.85"
$ws.Range("E26").Value = "  -1.86%  "
$ws.Range("E27").Value = "  -1.15%  "
$ws.Range("E28").Value = "  +0.08%  "
$ws.Range("D29").Value = "15.51"
$ws.Range("E29").Value = "  -1.07%  "
$ws.Range("E30").Value = "  -1.33%  "
$ws.Range("E31").Value = "  -0.90%  "
$ws.Range("E32").Value = "  -1.16%  "
$ws.Range("D33").Value = "1.466.86"
$ws.Range("E33").Value = "  +1.37%  "
$ws.Range("D34").Value = "3.06"
$ws.Range("E34").Value = "  -3.32%  "
$ws.Range("E35").Value = "  -4.16%  "
$ws.Range("E36").Value = "  -0.30%  "
$ws.Range("D37").Value = "0.950"
$ws.Range("E37").Value = "  +5.97%  "
$ws.Range("D38").Value = "0.557"
$ws.Range("E38").Value = "  -2.60%  "
$ws.Range("E39").Value = "  -0.84%  "
$ws.Range("D40").Value = "0.857"
$ws.Range("E40").Value = "  -3.22%  "
$ws.Range("E41").Value = "  +0.09%  "
$ws.Range("D42").Value = "67.87"
$ws.Range("E42").Value = "  +2.58%  "
$ws.Range("B43").Value = "mCoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin"
$ws.Range("D43").Value = "2.47"
$ws.Range("E43").Value = "  +0.66%  "
$ws.Range("B44").Value = "WEMIXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").Value = "0.987"
$ws.Range("E44").Value = "  -4.58%  "
$ws.Range("B45").Value = "MXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D45").Value = "2.20"
$ws.Range("E45").Value = "  -2.45%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").Value = "5.24"
$ws.Range("E46").Value = "  -7.82%  "
$ws.Range("B47").Value = "RocketPoolETH"
$ws.Range("C47").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D47").Value = "1.755.17"
$ws.Range("E47").Value = "  -1.81%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "1.72"
$ws.Range("E48").Value = "  +0.61%  "
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").Value = "86.42"
$ws.Range("E49").Value = "  -0.15%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₆0105"
$ws.Range("E50").Value = "  -2.65%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "0.101"
$ws.Range("E51").Value = "  +1.42%  "
